$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 676, shifting existing rows 676-687 down to 682-693
$ws.Rows.Item(676).Resize(6).Insert()

# Fill the newly inserted rows 676-681 with the new week of data
# Row 676
$ws.Cells.Item(676,1).Value = 9
$ws.Cells.Item(676,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(676,3).Value = 'Metropolitana'
$ws.Cells.Item(676,4).Value = 44628
$ws.Cells.Item(676,5).Value = 13
$ws.Cells.Item(676,6).Value = 100112027
$ws.Cells.Item(676,7).Value = 'Melón'
$ws.Cells.Item(676,8).Value = 'Calameño'
$ws.Cells.Item(676,9).Value = 'Extra'
$ws.Cells.Item(676,10).Value = 160
$ws.Cells.Item(676,11).Value = 1100
$ws.Cells.Item(676,12).Value = 1200
$ws.Cells.Item(676,13).Value = 1150
$ws.Cells.Item(676,14).Value = '$/unidad'
$ws.Cells.Item(676,15).Value = 'Región Metropolitana'
$ws.Cells.Item(676,16).Value = 1150
$ws.Cells.Item(676,17).Value = 1
$ws.Cells.Item(676,18).Value = 'Hortaliza'

# Row 677
$ws.Cells.Item(677,1).Value = 9
$ws.Cells.Item(677,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(677,3).Value = 'Metropolitana'
$ws.Cells.Item(677,4).Value = 44628
$ws.Cells.Item(677,5).Value = 13
$ws.Cells.Item(677,6).Value = 100112027
$ws.Cells.Item(677,7).Value = 'Melón'
$ws.Cells.Item(677,8).Value = 'Calameño'
$ws.Cells.Item(677,9).Value = 'Primera'
$ws.Cells.Item(677,10).Value = 250
$ws.Cells.Item(677,11).Value = 900
$ws.Cells.Item(677,12).Value = 1000
$ws.Cells.Item(677,13).Value = 950
$ws.Cells.Item(677,14).Value = '$/unidad'
$ws.Cells.Item(677,15).Value = 'Región Metropolitana'
$ws.Cells.Item(677,16).Value = 950
$ws.Cells.Item(677,17).Value = 1
$ws.Cells.Item(677,18).Value = 'Hortaliza'

# Row 678
$ws.Cells.Item(678,1).Value = 9
$ws.Cells.Item(678,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(678,3).Value = 'Metropolitana'
$ws.Cells.Item(678,4).Value = 44628
$ws.Cells.Item(678,5).Value = 13
$ws.Cells.Item(678,6).Value = 100112027
$ws.Cells.Item(678,7).Value = 'Melón'
$ws.Cells.Item(678,8).Value = 'Calameño'
$ws.Cells.Item(678,9).Value = 'Segunda'
$ws.Cells.Item(678,10).Value = 97
$ws.Cells.Item(678,11).Value = 700
$ws.Cells.Item(678,12).Value = 800
$ws.Cells.Item(678,13).Value = 749
$ws.Cells.Item(678,14).Value = '$/unidad'
$ws.Cells.Item(678,15).Value = 'Región Metropolitana'
$ws.Cells.Item(678,16).Value = 749
$ws.Cells.Item(678,17).Value = 1
$ws.Cells.Item(678,18).Value = 'Hortaliza'

# Row 679
$ws.Cells.Item(679,1).Value = 9
$ws.Cells.Item(679,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(679,3).Value = 'Metropolitana'
$ws.Cells.Item(679,4).Value = 44628
$ws.Cells.Item(679,5).Value = 13
$ws.Cells.Item(679,6).Value = 100112027
$ws.Cells.Item(679,7).Value = 'Melón'
$ws.Cells.Item(679,8).Value = 'Tuna'
$ws.Cells.Item(679,9).Value = 'Extra'
$ws.Cells.Item(679,10).Value = 160
$ws.Cells.Item(679,11).Value = 1100
$ws.Cells.Item(679,12).Value = 1200
$ws.Cells.Item(679,13).Value = 1150
$ws.Cells.Item(679,14).Value = '$/unidad'
$ws.Cells.Item(679,15).Value = 'Región Metropolitana'
$ws.Cells.Item(679,16).Value = 1150
$ws.Cells.Item(679,17).Value = 1
$ws.Cells.Item(679,18).Value = 'Hortaliza'

# Row 680
$ws.Cells.Item(680,1).Value = 9
$ws.Cells.Item(680,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(680,3).Value = 'Metropolitana'
$ws.Cells.Item(680,4).Value = 44628
$ws.Cells.Item(680,5).Value = 13
$ws.Cells.Item(680,6).Value = 100112027
$ws.Cells.Item(680,7).Value = 'Melón'
$ws.Cells.Item(680,8).Value = 'Tuna'
$ws.Cells.Item(680,9).Value = 'Primera'
$ws.Cells.Item(680,10).Value = 340
$ws.Cells.Item(680,11).Value = 900
$ws.Cells.Item(680,12).Value = 1000
$ws.Cells.Item(680,13).Value = 950
$ws.Cells.Item(680,14).Value = '$/unidad'
$ws.Cells.Item(680,15).Value = 'Región Metropolitana'
$ws.Cells.Item(680,16).Value = 950
$ws.Cells.Item(680,17).Value = 1
$ws.Cells.Item(680,18).Value = 'Hortaliza'

# Row 681
$ws.Cells.Item(681,1).Value = 9
$ws.Cells.Item(681,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(681,3).Value = 'Metropolitana'
$ws.Cells.Item(681,4).Value = 44628
$ws.Cells.Item(681,5).Value = 13
$ws.Cells.Item(681,6).Value = 100112027
$ws.Cells.Item(681,7).Value = 'Melón'
$ws.Cells.Item(681,8).Value = 'Tuna'
$ws.Cells.Item(681,9).Value = 'Segunda'
$ws.Cells.Item(681,10).Value = 250
$ws.Cells.Item(681,11).Value = 700
$ws.Cells.Item(681,12).Value = 800
$ws.Cells.Item(681,13).Value = 750
$ws.Cells.Item(681,14).Value = '$/unidad'
$ws.Cells.Item(681,15).Value = 'Región Metropolitana'
$ws.Cells.Item(681,16).Value = 750
$ws.Cells.Item(681,17).Value = 1
$ws.Cells.Item(681,18).Value = 'Hortaliza'
